{"js": "// Lit Review edit: add a new, empty, justified paragraph right after the\n// paragraph discussing the dose/image-quality tradeoff (the paragraph that\n// ends with \"...justifies the need to improve the modelling of low dose CT\n// images.\"), inserting it before the existing blank paragraph that follows.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"justifies the need to improve the modelling of low dose CT images.\";\n\nlet targetParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.indexOf(marker) !== -1) {\n    targetParagraph = paragraph;\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the dose/image-quality tradeoff paragraph.\");\n}\n\n// Insert a brand-new empty paragraph right after the target paragraph.\nconst newParagraph = targetParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nnewParagraph.alignment = Word.Alignment.justified;\n\nawait context.sync();\n", "ps1": "# Lit Review edit: add a new, empty, justified paragraph right after the\n# paragraph discussing the dose/image-quality tradeoff (the paragraph that\n# ends with \"...justifies the need to improve the modelling of low dose CT\n# images.\"), inserting it before the existing blank paragraph that follows.\n\n$d = $word.ActiveDocument\n\n$marker = \"justifies the need to improve the modelling of low dose CT images.\"\n\n$targetParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*$marker*\") {\n        $targetParagraph = $p\n        break\n    }\n}\n\nif ($targetParagraph -eq $null) {\n    throw \"Could not locate the dose/image-quality tradeoff paragraph.\"\n}\n\n$targetIndex = $targetParagraph.Index\n\n# Insert a brand-new empty paragraph right after the target paragraph.\n$targetParagraph.Range.InsertParagraphAfter()\n\n# wdAlignParagraphJustify = 3\n$newParagraph = $d.Paragraphs.Item($targetIndex + 1)\n$newParagraph.Alignment = 3\n"}
